$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert new row 2: ZA7848 / 96.3 / January-February 2022 / COVID-19 Pandemic
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "ZA7848"
$ws.Range("B2").Value = "'96.3"
$ws.Range("C2").Value = "January-February 2022"
$ws.Range("D2").Value = "COVID-19 Pandemic"

# Insert new row 9: ZA7749 / 94.1 / October-November 2020 / Future of Europe, Democracy in the EU, and Values and Identities of EU citizens
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "ZA7749"
$ws.Range("B9").Value = "'94.1"
$ws.Range("C9").Value = "October-November 2020"
$ws.Range("D9").Value = "Future of Europe, Democracy in the EU, and Values and Identities of EU citizens"

# Match the author's final selection state
[void]$ws.Range("A10:B10").Select()
